$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert 3 new columns (G,H,I) and shift old Totaal/Winnaar ---
$ws.Range('F1').Value = 'Totaal Score'
$ws.Range('G1').Value = 'Aantal Darts'
$ws.Range('H1').Value = '3-Darts Gemiddelde'
$ws.Range('I1').Value = 'Totaal'
$ws.Range('J1').Value = 'Winnaar'

# Copy the header formatting (bold font, borders, centered alignment) from an
# existing styled header cell (A1) onto the three newly-introduced header cells.
$ws.Range('A1').Copy()
$ws.Range('H1:J1').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (rows 2-45): full refresh of Score/180s/100+/Totaal Score/Aantal
# Darts/3-Darts Gemiddelde/Totaal/Winnaar per the updated totaalstand numbers ---
$r = 2
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 'Nick Fitzpatrick'
$ws.Cells.Item($r, 3).Value = 20
$ws.Cells.Item($r, 4).Value = 4
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 9517
$ws.Cells.Item($r, 7).Value = 466
$ws.Cells.Item($r, 8).Value = 61.27
$ws.Cells.Item($r, 9).Value = 24
$ws.Cells.Item($r, 10).Value = 1

$r = 3
$ws.Cells.Item($r, 1).Value = 2
$ws.Cells.Item($r, 2).Value = 'Patrick La Gordt Dillié'
$ws.Cells.Item($r, 3).Value = 15
$ws.Cells.Item($r, 4).Value = 1
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 11764
$ws.Cells.Item($r, 7).Value = 553
$ws.Cells.Item($r, 8).Value = 63.82
$ws.Cells.Item($r, 9).Value = 16
$ws.Cells.Item($r, 10).Value = 0

$r = 4
$ws.Cells.Item($r, 1).Value = 3
$ws.Cells.Item($r, 2).Value = 'Jente Stienstra'
$ws.Cells.Item($r, 3).Value = 11
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 7891
$ws.Cells.Item($r, 7).Value = 438
$ws.Cells.Item($r, 8).Value = 54.05
$ws.Cells.Item($r, 9).Value = 11
$ws.Cells.Item($r, 10).Value = 0

$r = 5
$ws.Cells.Item($r, 1).Value = 3
$ws.Cells.Item($r, 2).Value = 'Robin Hood'
$ws.Cells.Item($r, 3).Value = 11
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 3733
$ws.Cells.Item($r, 7).Value = 214
$ws.Cells.Item($r, 8).Value = 52.33
$ws.Cells.Item($r, 9).Value = 11
$ws.Cells.Item($r, 10).Value = 0

$r = 6
$ws.Cells.Item($r, 1).Value = 3
$ws.Cells.Item($r, 2).Value = 'Chris C'
$ws.Cells.Item($r, 3).Value = 9
$ws.Cells.Item($r, 4).Value = 1
$ws.Cells.Item($r, 5).Value = 1
$ws.Cells.Item($r, 6).Value = 6422
$ws.Cells.Item($r, 7).Value = 314
$ws.Cells.Item($r, 8).Value = 61.36
$ws.Cells.Item($r, 9).Value = 11
$ws.Cells.Item($r, 10).Value = 0

$r = 7
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = 'Mark Riches'
$ws.Cells.Item($r, 3).Value = 10
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 6838
$ws.Cells.Item($r, 7).Value = 375
$ws.Cells.Item($r, 8).Value = 54.7
$ws.Cells.Item($r, 9).Value = 10
$ws.Cells.Item($r, 10).Value = 0

$r = 8
$ws.Cells.Item($r, 1).Value = 7
$ws.Cells.Item($r, 2).Value = 'luca Simon'
$ws.Cells.Item($r, 3).Value = 8
$ws.Cells.Item($r, 4).Value = 1
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 5455
$ws.Cells.Item($r, 7).Value = 294
$ws.Cells.Item($r, 8).Value = 55.66
$ws.Cells.Item($r, 9).Value = 9
$ws.Cells.Item($r, 10).Value = 0

$r = 9
$ws.Cells.Item($r, 1).Value = 8
$ws.Cells.Item($r, 2).Value = 'Powy 🏴󠁧󠁢󠁷󠁬󠁳󠁿'
$ws.Cells.Item($r, 3).Value = 7
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 3013
$ws.Cells.Item($r, 7).Value = 146
$ws.Cells.Item($r, 8).Value = 61.91
$ws.Cells.Item($r, 9).Value = 7
$ws.Cells.Item($r, 10).Value = 0

$r = 10
$ws.Cells.Item($r, 1).Value = 8
$ws.Cells.Item($r, 2).Value = 'Oscar Ebbeling'
$ws.Cells.Item($r, 3).Value = 6
$ws.Cells.Item($r, 4).Value = 1
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 3800
$ws.Cells.Item($r, 7).Value = 237
$ws.Cells.Item($r, 8).Value = 48.1
$ws.Cells.Item($r, 9).Value = 7
$ws.Cells.Item($r, 10).Value = 0

$r = 11
$ws.Cells.Item($r, 1).Value = 10
$ws.Cells.Item($r, 2).Value = 'Gijs Tromp'
$ws.Cells.Item($r, 3).Value = 6
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 4553
$ws.Cells.Item($r, 7).Value = 232
$ws.Cells.Item($r, 8).Value = 58.88
$ws.Cells.Item($r, 9).Value = 6
$ws.Cells.Item($r, 10).Value = 0

$r = 12
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 'Dávid Márfi'
$ws.Cells.Item($r, 3).Value = 5
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 2820
$ws.Cells.Item($r, 7).Value = 186
$ws.Cells.Item($r, 8).Value = 45.48
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 10).Value = 0

$r = 13
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 'Matthew Cooke'
$ws.Cells.Item($r, 3).Value = 5
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 3112
$ws.Cells.Item($r, 7).Value = 231
$ws.Cells.Item($r, 8).Value = 40.42
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 10).Value = 0

$r = 14
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 'Pascal Ritter'
$ws.Cells.Item($r, 3).Value = 5
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1865
$ws.Cells.Item($r, 7).Value = 116
$ws.Cells.Item($r, 8).Value = 48.23
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 10).Value = 0

$r = 15
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 'Ayden Veenstra'
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 1
$ws.Cells.Item($r, 6).Value = 3423
$ws.Cells.Item($r, 7).Value = 176
$ws.Cells.Item($r, 8).Value = 58.35
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 10).Value = 0

$r = 16
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 'Chris Leo'
$ws.Cells.Item($r, 3).Value = 3
$ws.Cells.Item($r, 4).Value = 1
$ws.Cells.Item($r, 5).Value = 1
$ws.Cells.Item($r, 6).Value = 2261
$ws.Cells.Item($r, 7).Value = 130
$ws.Cells.Item($r, 8).Value = 52.18
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 10).Value = 0

$r = 17
$ws.Cells.Item($r, 1).Value = 16
$ws.Cells.Item($r, 2).Value = 'Noah B'
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 3297
$ws.Cells.Item($r, 7).Value = 219
$ws.Cells.Item($r, 8).Value = 45.16
$ws.Cells.Item($r, 9).Value = 4
$ws.Cells.Item($r, 10).Value = 0

$r = 18
$ws.Cells.Item($r, 1).Value = 16
$ws.Cells.Item($r, 2).Value = 'Stephen D'
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1710
$ws.Cells.Item($r, 7).Value = 102
$ws.Cells.Item($r, 8).Value = 50.29
$ws.Cells.Item($r, 9).Value = 4
$ws.Cells.Item($r, 10).Value = 0

$r = 19
$ws.Cells.Item($r, 1).Value = 18
$ws.Cells.Item($r, 2).Value = 'Flo Chételat'
$ws.Cells.Item($r, 3).Value = 3
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1420
$ws.Cells.Item($r, 7).Value = 104
$ws.Cells.Item($r, 8).Value = 40.96
$ws.Cells.Item($r, 9).Value = 3
$ws.Cells.Item($r, 10).Value = 0

$r = 20
$ws.Cells.Item($r, 1).Value = 18
$ws.Cells.Item($r, 2).Value = 'Joshua Taylor'
$ws.Cells.Item($r, 3).Value = 3
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1494
$ws.Cells.Item($r, 7).Value = 112
$ws.Cells.Item($r, 8).Value = 40.02
$ws.Cells.Item($r, 9).Value = 3
$ws.Cells.Item($r, 10).Value = 0

$r = 21
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Alexandra Kerr'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1778
$ws.Cells.Item($r, 7).Value = 109
$ws.Cells.Item($r, 8).Value = 48.94
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 22
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Angelo M'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 2238
$ws.Cells.Item($r, 7).Value = 120
$ws.Cells.Item($r, 8).Value = 55.95
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 23
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Blind Eagle'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = ''
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 24
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Daniel Maddison'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1765
$ws.Cells.Item($r, 7).Value = 121
$ws.Cells.Item($r, 8).Value = 43.76
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 25
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Dennis Sangler'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 822
$ws.Cells.Item($r, 7).Value = 54
$ws.Cells.Item($r, 8).Value = 45.67
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 26
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Evan Keeping'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1721
$ws.Cells.Item($r, 7).Value = 106
$ws.Cells.Item($r, 8).Value = 48.71
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 27
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Jakub Vraspír'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 882
$ws.Cells.Item($r, 7).Value = 54
$ws.Cells.Item($r, 8).Value = 49
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 28
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Jordan M'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 423
$ws.Cells.Item($r, 7).Value = 27
$ws.Cells.Item($r, 8).Value = 47
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 29
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Jun Denila'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 944
$ws.Cells.Item($r, 7).Value = 57
$ws.Cells.Item($r, 8).Value = 49.68
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 30
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Luke Kelly'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = ''
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 31
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Manolito Verleyen'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 978
$ws.Cells.Item($r, 7).Value = 75
$ws.Cells.Item($r, 8).Value = 39.12
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 32
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Nathan May'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 2457
$ws.Cells.Item($r, 7).Value = 124
$ws.Cells.Item($r, 8).Value = 59.44
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 33
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 'Robbie Hennes'
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 862
$ws.Cells.Item($r, 7).Value = 63
$ws.Cells.Item($r, 8).Value = 41.05
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 0

$r = 34
$ws.Cells.Item($r, 1).Value = 33
$ws.Cells.Item($r, 2).Value = 'Shea Bannon'
$ws.Cells.Item($r, 3).Value = 1
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1360
$ws.Cells.Item($r, 7).Value = 95
$ws.Cells.Item($r, 8).Value = 42.95
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = 0

$r = 35
$ws.Cells.Item($r, 1).Value = 33
$ws.Cells.Item($r, 2).Value = 'Yannick den Daggelder'
$ws.Cells.Item($r, 3).Value = 1
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 1265
$ws.Cells.Item($r, 7).Value = 78
$ws.Cells.Item($r, 8).Value = 48.65
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = 0

$r = 36
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Albin L'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 808
$ws.Cells.Item($r, 7).Value = 42
$ws.Cells.Item($r, 8).Value = 57.71
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 37
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Alfie Martin'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 888
$ws.Cells.Item($r, 7).Value = 51
$ws.Cells.Item($r, 8).Value = 52.24
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 38
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Bernhard Dierschke'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = ''
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 39
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Diego Meerveld'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 929
$ws.Cells.Item($r, 7).Value = 72
$ws.Cells.Item($r, 8).Value = 38.71
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 40
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Jaiden Powell'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 779
$ws.Cells.Item($r, 7).Value = 42
$ws.Cells.Item($r, 8).Value = 55.64
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 41
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Jay De-Winton'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 976
$ws.Cells.Item($r, 7).Value = 69
$ws.Cells.Item($r, 8).Value = 42.43
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 42
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Martin Koch'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 858
$ws.Cells.Item($r, 7).Value = 60
$ws.Cells.Item($r, 8).Value = 42.9
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 43
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Quinten the Bigfoot Dijkstra'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 0
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = ''
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 44
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Tom Jones'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 924
$ws.Cells.Item($r, 7).Value = 54
$ws.Cells.Item($r, 8).Value = 51.33
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

$r = 45
$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = 'Tomm Gordon'
$ws.Cells.Item($r, 3).Value = 0
$ws.Cells.Item($r, 4).Value = 0
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 789
$ws.Cells.Item($r, 7).Value = 51
$ws.Cells.Item($r, 8).Value = 46.41
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0

